# Auto-generated Excel COM-interop script applying scheduled price-refresh updates
# to the per-job (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) profit tables.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 634.63635
$ws.Range("I2").Value = 396.2
$ws.Range("K2").Value = 396.2
$ws.Range("M2").Value = -283.2
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H132").Value = 39684956
$ws.Range("I132").Value = 8548377
$ws.Range("J132").Value = 66669988
$ws.Range("K132").Value = 25645131
$ws.Range("L132").Value = 200009964
$ws.Range("M132").Value = -25642601
$ws.Range("N132").Value = -200015024
$ws.Range("H137").Value = 618651.5600000001
$ws.Range("I137").Value = 1734.72
$ws.Range("J137").Value = 855927.25
$ws.Range("K137").Value = 5204.16
$ws.Range("L137").Value = 2567781.75
$ws.Range("M137").Value = -2654.16
$ws.Range("N137").Value = -2572881.75

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H119").Value = 29685
$ws.Range("J119").Value = 29685
$ws.Range("L119").Value = 29685
$ws.Range("N119").Value = -39361
$ws.Range("H125").Value = 34988.89
$ws.Range("J125").Value = 34988.89
$ws.Range("L125").Value = 34988.89
$ws.Range("N125").Value = -44828.89
$ws.Range("H139").Value = 28904.334
$ws.Range("J139").Value = 28904.334
$ws.Range("L139").Value = 28904.334
$ws.Range("N139").Value = -39184.334

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 10000
$ws.Range("J16").Value = 10000
$ws.Range("L16").Value = 10000
$ws.Range("N16").Value = -10340
$ws.Range("H105").Value = 2530.182
$ws.Range("I105").Value = 2523.0952
$ws.Range("K105").Value = 2523.0952
$ws.Range("M105").Value = -776.0952000000002
$ws.Range("H138").Value = 49572.75
$ws.Range("J138").Value = 49572.75
$ws.Range("L138").Value = 49572.75
$ws.Range("N138").Value = -59852.75

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10871036
$ws.Range("I31").Value = 1029.5714
$ws.Range("J31").Value = 20001840
$ws.Range("K31").Value = 1029.5714
$ws.Range("L31").Value = 20001840
$ws.Range("M31").Value = -734.5714
$ws.Range("N31").Value = -20002430
$ws.Range("H34").Value = 10871036
$ws.Range("I34").Value = 1029.5714
$ws.Range("J34").Value = 20001840
$ws.Range("K34").Value = 1029.5714
$ws.Range("L34").Value = 20001840
$ws.Range("M34").Value = -827.5714
$ws.Range("N34").Value = -20002244
$ws.Range("H132").Value = 4296.4614
$ws.Range("I132").Value = 4149.1113
$ws.Range("K132").Value = 12447.3339
$ws.Range("M132").Value = -9917.333899999998
$ws.Range("H134").Value = 34376776
$ws.Range("I134").Value = 4349655.5
$ws.Range("J134").Value = 111112750
$ws.Range("K134").Value = 13048966.5
$ws.Range("L134").Value = 333338250
$ws.Range("M134").Value = -13046431.5
$ws.Range("N134").Value = -333343320

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 11828535
$ws.Range("I5").Value = 10000337
$ws.Range("J5").Value = 15152531
$ws.Range("K5").Value = 30001011
$ws.Range("L5").Value = 45457593
$ws.Range("M5").Value = -30000899
$ws.Range("N5").Value = -45457817
$ws.Range("H36").Value = 2150.0833
$ws.Range("I36").Value = 1942.5714
$ws.Range("J36").Value = 2440.6
$ws.Range("K36").Value = 5827.7142
$ws.Range("L36").Value = 7321.799999999999
$ws.Range("M36").Value = -5658.7142
$ws.Range("N36").Value = -7659.799999999999
$ws.Range("H39").Value = 3356.875
$ws.Range("J39").Value = 3722.1428
$ws.Range("L39").Value = 11166.4284
$ws.Range("N39").Value = -11754.4284
$ws.Range("H95").Value = 9500
$ws.Range("J95").Value = 9500
$ws.Range("L95").Value = 28500
$ws.Range("N95").Value = -32618
$ws.Range("H122").Value = 714.64514
$ws.Range("I122").Value = 478.33334
$ws.Range("J122").Value = 936.1875
$ws.Range("K122").Value = 4305.00006
$ws.Range("L122").Value = 8425.6875
$ws.Range("M122").Value = -1855.00006
$ws.Range("N122").Value = -13325.6875
$ws.Range("H131").Value = 891.12195
$ws.Range("J131").Value = 971.7361
$ws.Range("L131").Value = 2915.2083
$ws.Range("N131").Value = -12995.2083
$ws.Range("H135").Value = 11828535
$ws.Range("I135").Value = 10000337
$ws.Range("J135").Value = 15152531
$ws.Range("K135").Value = 90003033
$ws.Range("L135").Value = 136372779
$ws.Range("M135").Value = -90000498
$ws.Range("N135").Value = -136377849

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 23446
$ws.Range("J12").Value = 6889
$ws.Range("L12").Value = 6889
$ws.Range("N12").Value = -7169
$ws.Range("H126").Value = 1950
$ws.Range("I126").Value = 2200
$ws.Range("J126").Value = 1700
$ws.Range("K126").Value = 6600
$ws.Range("L126").Value = 5100
$ws.Range("M126").Value = -4130
$ws.Range("N126").Value = -10040
$ws.Range("H135").Value = 37500
$ws.Range("J135").Value = 37500
$ws.Range("L135").Value = 37500
$ws.Range("N135").Value = -47640

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1778.1111
$ws.Range("I16").Value = 1625.125
$ws.Range("J16").Value = 3002
$ws.Range("K16").Value = 1625.125
$ws.Range("L16").Value = 3002
$ws.Range("M16").Value = -1455.125
$ws.Range("N16").Value = -3342
$ws.Range("H136").Value = 1401.8392
$ws.Range("I136").Value = 1383.122
$ws.Range("J136").Value = 1453
$ws.Range("K136").Value = 4149.366
$ws.Range("L136").Value = 4359
$ws.Range("M136").Value = -1599.366
$ws.Range("N136").Value = -9459

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H29").Value = 10402.5
$ws.Range("J29").Value = 10866.667
$ws.Range("L29").Value = 10866.667
$ws.Range("N29").Value = -11446.667
$ws.Range("H32").Value = 8139.4
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 8139.4
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 8139.4
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -8773.4
$ws.Range("H33").Value = 17999.666
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 17999.666
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 17999.666
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -18499.666
$ws.Range("H36").Value = 17999.666
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 17999.666
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 17999.666
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -18499.666
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H107").Value = 9719.272000000001
$ws.Range("I107").Value = 11711.777
$ws.Range("K107").Value = 35135.331
$ws.Range("M107").Value = -33215.331
$ws.Range("H132").Value = 3692.6758
$ws.Range("I132").Value = 4578.143
$ws.Range("J132").Value = 2530.5
$ws.Range("K132").Value = 13734.429
$ws.Range("L132").Value = 7591.5
$ws.Range("M132").Value = -11204.429
$ws.Range("N132").Value = -12651.5

